$d = $word.ActiveDocument

# --- First paragraph: update paragraph formatting ---
$para = $d.Paragraphs(1)

# Add paragraph border (top/left/bottom/right) with 5pt text-distance spacing,
# matching the <w:pBdr><w:top w:space="5"/> ... structure.
$para.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$para.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$para.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$para.Range.ParagraphFormat.Borders.DistanceFromRight = 5

# Change left indent from 120 twips to 225 twips (1 point = 20 twips).
$para.Range.ParagraphFormat.LeftIndent = 225 / 20

# --- Update the bookmark/ID text and drop the trailing space run ---
$rng = $d.Content
$found = $rng.Find.Execute("**ID__AFFARS_mp_5315_3_topic_5__ID**", $true, $false, $false, $false, $false,
                            $true, 1, $false)
$rng.Text = "**ID__AFFARS_MP_5315_3_1_3__ID**"

# The old text was immediately followed by a run containing a single space;
# remove that now-orphaned trailing space character.
$spaceRng = $d.Range($rng.End, $rng.End + 1)
if ($spaceRng.Text -eq " ") {
    $spaceRng.Delete()
}
